$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$ws.Range("D2").Value = "37.366.06"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.072.80"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "234.15"
$ws.Range("E5").Value = "  -1.61%  "
Set-TextValue "D6" "0.626"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue "D8" "57.07"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D12").Value = "2.376.83"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  +1.06%  "
Set-TextValue "D14" "20.80"
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").Value = "2.075.38"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "37.296.97"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("E19").Value = "  +2.36%  "
Set-TextValue "D20" "69.47"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").Value = "0.0₃0815"
$ws.Range("E21").Value = "  -0.40%  "
Set-TextValue "D22" "227.16"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -3.42%  "
Set-TextValue "D26" "167.04"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  +3.13%  "
Set-TextValue "D29" "19.12"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("E31").Value = "  -0.94%  "
Set-TextValue "D32" "4.48"
$ws.Range("E32").Value = "  -0.61%  "
Set-TextValue "D33" "0.0618"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("E35").Value = "  -4.82%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("E38").Value = "  -4.18%  "
Set-TextValue "D39" "5.68"
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("E40").Value = "  -0.28%  "
Set-TextValue "D41" "4.40"
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("D42").Value = "1.478.82"
$ws.Range("E42").Value = "  -0.64%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "96.45"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D44" "0.0942"
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D45" "0.0212"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D46" "1.17"
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("E47").Value = "  -1.16%  "
Set-TextValue "D48" "15.05"
$ws.Range("E48").Value = "  -8.96%  "
$ws.Range("E49").Value = "  -3.30%  "
Set-TextValue "D50" "2.96"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "2.266.31"
$ws.Range("E51").Value = "  -0.59%  "
